# Apply the "add complex and fix input" edit:
#  - Rewrite the 5 question slides (1,3,5,7,9) with new riddles, dropping
#    the trailing "Класс: 1-4 класс" paragraph (and the stray blank one on
#    slide 1) so each question is a single paragraph.
#  - Rewrite the 5 matching answer slides (2,4,6,8,10) with new answers.
#  - Append two brand-new slides (11/12) with a new question/answer pair.

$p = $ppt.ActivePresentation

# ---- existing question slides -------------------------------------------------
$p.Slides.Item(1).Shapes.Item(2).TextFrame.TextRange.Text = "Вопрос: В настоящее время ИМ пользуется более 2,5 млрд человек. В 1991 году ОН стал общедоступным, сейчас есть почти в каждом доме. По статистике в современном мире каждая восьмая супружеская пара познакомилась в НЕМ. Назовите ЕГО."

$p.Slides.Item(3).Shapes.Item(2).TextFrame.TextRange.Text = "Вопрос: Самым первым материалом для создания ЭТОГО служил рыбий скелет. А впервые в привычном для нас виде ЭТО появилось в Древнем Риме и было сделано из слоновой кости. На Руси ЭТО делали из дерева. А сейчас ЭТО все чаще делают из пластика. Причем чаще всего ЭТИМ пользуются представительницы женского пола. Что ЭТО?"

$p.Slides.Item(5).Shapes.Item(2).TextFrame.TextRange.Text = "Вопрос: Это животное - единственное, которое не умеет зевать. Длина ЕГО языка достигает 50 см, а длина хвоста - 2,5 м. Удивительно то, что строение шеи этого животного не позволяет ему дотянуться до земли, поэтому ему приходится вставать на колени, чтобы поднять что-то с пола. Назовите это животное."

$p.Slides.Item(7).Shapes.Item(2).TextFrame.TextRange.Text = "Вопрос: В 1765 г. Екатерина II озаботилась помощью голодающим крестьянам Финляндии. После некоторых поисков, коллегия, которой был поручен этот вопрос, решила, что лучше всего использовать “земляные яблоки”. Что ЭТО было?"

$p.Slides.Item(9).Shapes.Item(2).TextFrame.TextRange.Text = "Вопрос: Существует множество легенд об этом явлении. Ему приписывают чудодейственные свойства. Говорят, что им можно исцелиться или избавиться от злых чар. Очень часто ЭТО называют “водяные алмазы”. Что ЭТО?"

# ---- existing answer slides -----------------------------------------------------
$p.Slides.Item(2).Shapes.Item(2).TextFrame.TextRange.Text = "Ответ: интернет"
$p.Slides.Item(4).Shapes.Item(2).TextFrame.TextRange.Text = "Ответ: расческа / гребень для волос"
$p.Slides.Item(6).Shapes.Item(2).TextFrame.TextRange.Text = "Ответ: жираф"
$p.Slides.Item(8).Shapes.Item(2).TextFrame.TextRange.Text = "Ответ: картофель / картошка"
$p.Slides.Item(10).Shapes.Item(2).TextFrame.TextRange.Text = "Ответ: роса"

# ---- two brand-new slides appended at the end ------------------------------------
$s11 = $p.Slides.Add(11, 2)
$s11.Shapes.Item(1).TextFrame.TextRange.Text = "ЧГК"
$s11.Shapes.Item(2).TextFrame.TextRange.Text = "Вопрос: Некоторые ученые в Древней Греции воспринимали натуральные числа как собрание ИХ. При этом саму ЕЕ числом они не считали. Отсюда пошла теория, что ОНА — это не начало числового ряда, а лишь некий начальный процесс. Такого мнения придерживался Платон. Назовите ЕЕ."

$s12 = $p.Slides.Add(12, 2)
$s12.Shapes.Item(2).TextFrame.TextRange.Text = "Ответ: единица"
